$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly record row at position 235 (pushes existing rows 235-286 down to 236-287)
$ws.Rows.Item(235).Insert()

$ws.Cells.Item(235, 1).Value = 4
$ws.Cells.Item(235, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(235, 3).Value = "Los Lagos"
$ws.Cells.Item(235, 4).Value = 44754
$ws.Cells.Item(235, 5).Value = 10
$ws.Cells.Item(235, 6).Value = 100112003
$ws.Cells.Item(235, 7).Value = "Ajo"
$ws.Cells.Item(235, 8).Value = "Chino"
$ws.Cells.Item(235, 9).Value = "Primera"
$ws.Cells.Item(235, 10).Value = 200
$ws.Cells.Item(235, 11).Value = 27000
$ws.Cells.Item(235, 12).Value = 28000
$ws.Cells.Item(235, 13).Value = 27500
$ws.Cells.Item(235, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(235, 15).Value = "China"
$ws.Cells.Item(235, 16).Value = 2750
$ws.Cells.Item(235, 17).Value = 10
$ws.Cells.Item(235, 18).Value = "Hortaliza"
